# edit.ps1 - Apply resume wording/content changes described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $old"
    }
}

# --- Header ---
Replace-Text "DAVID GLASS" "David Glass"
Replace-Text "Chicago, IL 60618 ⋄ 847.764.9200 ⋄ dglass2525@gmail.com ⋄ thedavidglass.com ⋄ github.com/dglass710 ⋄ linkedin.com/in/david-a-glass" "Chicago, IL 60618 ⋄ 847.764.9200 ⋄ dglass2525@gmail.com ⋄ thedavidglass.com ⋄ GitHub ⋄ LinkedIn"

# --- Objective ---
Replace-Text "To secure a position as an IT Systems Administrator, maintaining and optimizing enterprise systems and infrastructure while ensuring reliable and secure network operations across on-premises and cloud environments." "Seeking a cybersecurity analyst role within a dynamic Security Operations Center (SOC) to leverage skills in real-time threat monitoring, vulnerability scanning, and incident response."

# --- Certifications ---
Replace-Text "CompTIA Security+ Certified" "CompTIA Security+"

# --- Education: restructure into 4 paragraphs ---
# Original 3 paragraphs:
#   "Northwestern University Cybersecurity Program: Graduated June 2024"
#   "DePaul University: Graduated March 2022"
#   "    B.S. in Applied and Computational Mathematics; Minors in Computer Science and Physics."
# Target 4 paragraphs:
#   "Northwestern University"
#   "    Cybersecurity Program Certificate, Graduated June 2024"
#   "DePaul University"
#   "    Bachelor of Science in Applied and Computational Mathematics, Minors: Computer Science and Physics, Graduated March 2022"

# Step 1: rename the Northwestern line (drop the cert/date suffix)
Replace-Text "Northwestern University Cybersecurity Program: Graduated June 2024" "Northwestern University"

# Step 2: the old "DePaul University: Graduated March 2022" line becomes the new
# Northwestern certificate detail line
Replace-Text "DePaul University: Graduated March 2022" "    Cybersecurity Program Certificate, Graduated June 2024"

# Step 3: the old B.S. detail line becomes the new "DePaul University" line
Replace-Text "    B.S. in Applied and Computational Mathematics; Minors in Computer Science and Physics." "DePaul University"

# Step 4: insert a brand-new paragraph after the (now renamed) "DePaul University" line
# carrying the full degree detail text
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -eq "DePaul University`r") {
        $targetIdx = $idx
    }
}
if ($targetIdx -ge 0) {
    $p = $d.Paragraphs.Item($targetIdx)
    $rng = $p.Range
    $rng2 = $d.Range($rng.Start, $rng.End - 1)
    $rng2.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($targetIdx + 1)
    $newRng = $d.Range($newP.Range.Start, $newP.Range.End - 1)
    $newRng.Text = "    Bachelor of Science in Applied and Computational Mathematics, Minors: Computer Science and Physics, Graduated March 2022"
} else {
    Write-Host "WARNING: could not locate the renamed DePaul University paragraph"
}

# --- Core Competencies ---
Replace-Text "Linux Administration, macOS Administration, Windows Administration, Active Directory, Group Policy, System Hardening, User and Permissions Management, Firewalls (UFW, firewalld), Process Management, Backup and Archive Management, Cron Jobs, Scripting Maintenance Tasks, Cloud Infrastructure (Azure, AWS), Containerization, Load Balancers, SQL Query Design, Docker, Automation Scripting, Technical Communication, Problem Solving." "Risk Assessment, Threat Mitigation, Incident Response, Vulnerability Scanning, Penetration Testing, Network Security, Splunk, Security Onion, Log Analysis, Snort Rules, Ethical Hacking, Privilege Escalation, Exploit Development, Packet Analysis, SIEM Analysis, Incident Reporting, Firewall Policy Development, OSINT Techniques, Recon-ng, Advanced Nmap Scans, Digital Evidence Preservation, Hashcat, Google Dorking, Shodan."

# --- Professional Experience: Mathnasium ---
Replace-Text "Mathnasium – Mathematics Instructor (April 2023 – Present, Chicago, IL)" "Mathematics Instructor – Mathnasium (April 2023 – Present, Chicago, IL)"
Replace-Text "Instructed 370 students in mathematical concepts." "Taught 370 students mathematical concepts."
Replace-Text "Enhanced security of iPads used in instruction by implementing guided access controls." "Enhanced security for instructional devices by implementing guided access controls on iPads."

# --- Professional Experience: DePaul University Math Department ---
Replace-Text "DePaul University Math Department – Undergraduate Student Researcher (November 2020 – May 2022, Chicago, IL)" "Undergraduate Student Researcher – DePaul University Math Department (November 2020 – May 2022, Chicago, IL)"
Replace-Text "Developed Python tools for research, focusing on preventing rounding errors in fraction representations." "Developed Python tools to address rounding errors in fractional computations."
Replace-Text "Advanced the understanding of the Frobenius coin problem and computed symmetry in large data sets." "Researched advanced aspects of the Frobenius coin problem and symmetry in large datasets."

# --- Technical Projects ---
Replace-Text "Developed and implemented firewall policies using UFW and firewalld to ensure secure server configurations and restrict unauthorized access." "Conducted penetration tests and network traffic analysis using Metasploit, Nmap, and Wireshark to identify and exploit vulnerabilities."
Replace-Text "Automated user account creation and permissions assignment with Bash scripts, improving efficiency for onboarding processes.                " "Designed firewall policies with UFW and firewalld to secure network configurations."
Replace-Text "Built and managed a virtual private cloud (VPC) on Azure, deploying virtual machines and configuring secure remote access via SSH.            " "Configured and monitored Splunk SIEM to analyze security logs, detect anomalies, and mitigate potential threats."
Replace-Text "Deployed and maintained Docker containers for streamlined application development and environment consistency.                              " "Executed vulnerability scans in virtual environments, prioritizing remediation for high-risk exposures."
Replace-Text "Configured cron jobs to automate system maintenance, including backups, log rotation, and system updates.                                   " "Created a Dockerized offline version of the Have I Been Pwned database for secure local credential queries."

Write-Host "Edits applied."
